# Fruta / hortaliza, semanal
# Weekly refresh of the Sandia (watermelon) price records: several rows'
# data (date, variety, quality, volume, prices, trade unit, origin) are
# updated to reflect newly reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg
    )

    $ws.Cells.Item($Row, 4).Value = $Fecha        # D - Fecha
    $ws.Cells.Item($Row, 8).Value = $Variedad     # H - Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad      # I - Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen     # J - Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin   # K - Precio minimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMax   # L - Precio maximo
    $ws.Cells.Item($Row, 13).Value = $PrecioProm  # M - Precio promedio ponderado
    $ws.Cells.Item($Row, 14).Value = $Unidad      # N - Unidad de comercializacion
    $ws.Cells.Item($Row, 15).Value = $Origen      # O - Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg    # P - Precio $/Kg
}

Set-Row 3  44510 "Sin especificar"      "Primera" 250 800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
Set-Row 4  44223 "Americana O Klondike" "Extra"   340 2500 2500 2500 "`$/unidad"                      "Región de O'Higgins" 2500
Set-Row 5  44223 "Americana O Klondike" "Primera" 400 2000 2000 2000 "`$/unidad"                      "Región de O'Higgins" 2000
Set-Row 6  44223 "Americana O Klondike" "Segunda" 300 1500 1500 1500 "`$/unidad"                      "Región de O'Higgins" 1500
Set-Row 7  44223 "Americana O Klondike" "Tercera" 160 1000 1000 1000 "`$/unidad"                      "Región de O'Higgins" 1000
Set-Row 8  44217 "Sin especificar"      "Extra"   400 2500 2500 2500 "`$/unidad"                      "Región de O'Higgins" 2500
Set-Row 9  44217 "Sin especificar"      "Primera" 280 2000 2000 2000 "`$/unidad"                      "Región de O'Higgins" 2000
Set-Row 10 44194 "Sin especificar"      "Extra"   120 3500 3500 3500 "`$/unidad"                      "Región de O'Higgins" 3500
Set-Row 11 44194 "Sin especificar"      "Primera" 200 3000 3000 3000 "`$/unidad"                      "Región de O'Higgins" 3000
Set-Row 12 44497 "Sin especificar"      "Primera" 250 800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
Set-Row 14 44167 "Sin especificar"      "Primera" 400 5000 5000 5000 "`$/unidad"                      "Región de O'Higgins" 5000
Set-Row 15 44167 "Sin especificar"      "Segunda" 560 3000 3000 3000 "`$/unidad"                      "Región de O'Higgins" 3000
Set-Row 16 44167 "Sin especificar"      "Tercera" 450 2000 2000 2000 "`$/unidad"                      "Región de O'Higgins" 2000
Set-Row 17 44495 "Sin especificar"      "Primera" 200 800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
Set-Row 18 44305 "Sin especificar"      "Primera" 100 2500 2500 2500 "`$/unidad"                      "Perú"                 2500
Set-Row 19 44477 "Sin especificar"      "Primera" 80  800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
Set-Row 20 44491 "Sin especificar"      "Primera" 150 800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
Set-Row 21 44488 "Sin especificar"      "Primera" 150 800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
Set-Row 22 44504 "Sin especificar"      "Primera" 200 800  800  800  "`$/kilo (volumen en unidades)" "Perú"                 800
